$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet (typo fix: "Demo Icons" -> "Demo Iconls")
$ws.Name = "Demo Iconls"

# Rename the four embedded picture shapes (cNvPr name regenerated)
$ws.Shapes.Item(1).Name = "0800271CCEE91EDC95A98A733BD1E059"
$ws.Shapes.Item(2).Name = "0800271CCEE91EDC95A98A733BD20059"
$ws.Shapes.Item(3).Name = "0800271CCEE91EDC95A98A733BD22059"
$ws.Shapes.Item(4).Name = "0800271CCEE91EDC95A98A733BD24059"
